$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap "Paises Bajos" / "Arabia Saudita" ordering (row 19 <-> row 20) ---
# and refresh their statistics at the same time.
$ws.Range("A19").Value = "Arabia Saudita"
$ws.Range("B19").Value = 44830
$ws.Range("C19").Value = 1905
$ws.Range("D19").Value = 17622
$ws.Range("E19").Value = 26935
$ws.Range("F19").Value = 147
$ws.Range("G19").Value = 9
$ws.Range("H19").Value = 273

$ws.Range("A20").Value = "Paises Bajos"
$ws.Range("B20").Value = 43211
$ws.Range("C20").Value = 227
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 37399
$ws.Range("F20").Value = 463
$ws.Range("G20").Value = 52
$ws.Range("H20").Value = 5562

# --- Swap "Singapur" / "Bielorrusia" ordering (row 29 <-> row 30) ---
# and refresh their statistics at the same time.
$ws.Range("A29").Value = "Bielorrusia"
$ws.Range("B29").Value = 25825
$ws.Range("C29").Value = 952
$ws.Range("D29").Value = 7711
$ws.Range("E29").Value = 17968
$ws.Range("F29").Value = 92
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 146

$ws.Range("A30").Value = "Singapur"
$ws.Range("B30").Value = 25346
$ws.Range("C30").Value = 675
$ws.Range("D30").Value = 3851
$ws.Range("E30").Value = 21474
$ws.Range("F30").Value = 20
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 21

# --- Plain statistic refreshes (no reordering) ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 1409769
$ws.Range("C4").Value = 1133
$ws.Range("E4").Value = 1029532
$ws.Range("G4").Value = 66
$ws.Range("H4").Value = 83491

# Alemania (row 11)
$ws.Range("B11").Value = 173524
$ws.Range("C11").Value = 353
$ws.Range("E11").Value = 17044
$ws.Range("G11").Value = 42
$ws.Range("H11").Value = 7780

# Suiza (row 25)
$ws.Range("E25").Value = 1743
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 1870

# Serbia (row 48)
$ws.Range("B48").Value = 10295
$ws.Range("C48").Value = 52
$ws.Range("D48").Value = 3824
$ws.Range("E48").Value = 6249
$ws.Range("F48").Value = 22
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 222

# Barein (row 59)
$ws.Range("B59").Value = 5780
$ws.Range("C59").Value = 249
$ws.Range("D59").Value = 2195
$ws.Range("E59").Value = 3576
$ws.Range("F59").Value = 7

# Guyana (row 162)
$ws.Range("D162").Value = 39
$ws.Range("E162").Value = 64

# --- Update the "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 15:05"
